$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from the last existing data row (28) down onto the new
# rows so the new cells pick up the same styles (centered numbers in
# column E, left-aligned text in columns F/G) without introducing any new
# cell-format (xf) entries in styles.xml.
$ws.Range("E28:G28").Copy()
$ws.Range("E29:G33").PasteSpecial(-4122)  # xlPasteFormats

# New "git branch" / merge-conflict related commands, continuing the table.
$ws.Cells.Item(29, 5).Value = 27
$ws.Cells.Item(30, 5).Value = 28
$ws.Cells.Item(31, 5).Value = 29
$ws.Cells.Item(32, 5).Value = 30
$ws.Cells.Item(33, 5).Value = 31

$ws.Cells.Item(29, 6).Value = "git branch -v"
$ws.Cells.Item(29, 7).Value = "This shows names of branch with their commit hash and commit message"

$ws.Cells.Item(30, 6).Value = "git branch --merged"
$ws.Cells.Item(31, 6).Value = "git branch --no-merged"

$ws.Cells.Item(30, 7).Value = "Already merged branches"
$ws.Cells.Item(31, 7).Value = "Not already merged branches"

$ws.Cells.Item(32, 6).Value = "git branch -d develop"
$ws.Cells.Item(32, 7).Value = "Gives error if develop is not merged"

$ws.Cells.Item(33, 6).Value = "git branch -D develop"
$ws.Cells.Item(33, 7).Value = "No error and branch gets deleted"

$ws.Range("G33").Select()
